function Set-RunText($rng, $text) {
    # Setting .Text alone on a sub-range tends to get folded back into
    # neighboring runs that share identical formatting once another edit
    # touches the same paragraph. Re-stamping the font properties (even to
    # their existing values) forces the engine to keep this text as its own
    # distinct run instead of merging it with its neighbors.
    $rng.Text = $text
    $rng.Font.Name = "Arial"
    $rng.Font.NameFarEast = "Microsoft Sans Serif"
    $rng.Font.NameBi = "Arial"
}

$d = $word.ActiveDocument
$p = $d.Paragraphs(8)

# ============================================================
# Phase 1 -- fix the actual wording (first person -> third person)
# ============================================================
$rng = $p.Range
$rng.Find.Execute("I enjoy") | Out-Null
$rng.Text = "Enjoys"

$rng = $p.Range
$rng.Find.Execute("taking my knowledge of IT to the next level. ") | Out-Null
$rng.Text = "taking knowledge of IT to the next level. "

$rng = $p.Range
$rng.Find.Execute("I can") | Out-Null
$rng.Text = "Can"

$rng = $p.Range
$rng.Find.Execute("I am detail-oriented") | Out-Null
$rng.Text = "Is detail-oriented"

$rng = $p.Range
$rng.Find.Execute("Finally, I am a valuable") | Out-Null
$rng.Text = "Finally, is a valuable"

# ============================================================
# Phase 2 -- re-split the paragraph into the same fine-grained runs the
# author's edit produced, and relocate the "_GoBack" bookmark.
# ============================================================

# "...and following up." -> "...and f" | [[_GoBack]] | "ollowing up."
$rngF = $p.Range
$rngF.Find.Execute("following up") | Out-Null
$collapsed = $rngF.Duplicate
$collapsed.Start = $rngF.Start + 1
$collapsed.End = $collapsed.Start
$d.Bookmarks.Add("_GoBack", $collapsed)

# "Enjoys"
$rng = $p.Range
$rng.Find.Execute("Enjoys") | Out-Null
Set-RunText $rng "Enjoys"

# "C"
$rng0 = $p.Range
$rng0.Find.Execute("Can disassemble") | Out-Null
$rng = $rng0.Duplicate
$rng.End = $rng.Start + 1
Set-RunText $rng "C"

# "Is"
$rng0 = $p.Range
$rng0.Find.Execute("Is detail-oriented") | Out-Null
$rng = $rng0.Duplicate
$rng.End = $rng.Start + 2
Set-RunText $rng "Is"

# " detail-oriented and always strive fo"
$rng = $p.Range
$rng.Find.Execute(" detail-oriented and always strive fo") | Out-Null
Set-RunText $rng " detail-oriented and always strive fo"

# "r perfection by "
$rng = $p.Range
$rng.Find.Execute("r perfection by ") | Out-Null
Set-RunText $rng "r perfection by "

# "utilizing"
$rng = $p.Range
$rng.Find.Execute("utilizing") | Out-Null
Set-RunText $rng "utilizing"

# " feedback"
$rng = $p.Range
$rng.Find.Execute(" feedback") | Out-Null
Set-RunText $rng " feedback"

# ". Finally, is"
$rng = $p.Range
$rng.Find.Execute(". Finally, is") | Out-Null
Set-RunText $rng ". Finally, is"

# " a valuable team member "
$rng = $p.Range
$rng.Find.Execute(" a valuable team member ") | Out-Null
Set-RunText $rng " a valuable team member "

# "who"
$rng = $p.Range
$rng.Find.Execute("who") | Out-Null
Set-RunText $rng "who"

# " is especially good at breaking the ice, bringing up new ideas, and f"
$rng = $p.Range
$rng.Find.Execute(" is especially good at breaking the ice, bringing up new ideas, and f") | Out-Null
Set-RunText $rng " is especially good at breaking the ice, bringing up new ideas, and f"

# ============================================================
# Merge the "Reimaging, Troubleshooting" / ", Automation, Scripting,
# Customer Service " runs into a single run (no wording change). Using
# Find's own Replace argument (rather than a manual Range.Text=
# assignment) reliably coalesces the two runs into one.
# ============================================================
$rngSkills = $d.Content
$rngSkills.Find.Execute("Reimaging, Troubleshooting", $true, $false, $false, $false, $false, $true, 1, $false, "Reimaging, Troubleshooting", 2) | Out-Null

Write-Output $p.Range.Text
